$d = $word.ActiveDocument
$replacements = @(
    ,@('2025-07-22 Tuesday', '2025-07-23 Wednesday')
    ,@('31+64=', '43-33=')
    ,@('44-7=', '37+18=')
    ,@('63+12=', '68-68=')
    ,@('87-52=', '14+50=')
    ,@('53+37=', '69+2=')
    ,@('84-8=', '35+44=')
    ,@('15+11=', '18+14=')
    ,@('83-13=', '90-9=')
    ,@('36-1=', '63-34=')
    ,@('78-22=', '78-35=')
    ,@('99-22=', '56-54=')
    ,@('85-56=', '49+27=')
    ,@('37-9=', '70-44=')
    ,@('59+7=', '27+28=')
    ,@('52+42=', '89-3=')
    ,@('53-7=', '68+24=')
    ,@('3+6=', '89-63=')
    ,@('14+23=', '53-52=')
    ,@('26-13=', '15+65=')
    ,@('48+35=', '37-35=')
    ,@('22+33=', '3+22=')
    ,@('40-28=', '42-17=')
    ,@('57+7=', '34+22=')
    ,@('51-37=', '25+39=')
    ,@('59-35=', '73-6=')
    ,@('41-40=', '47-9=')
    ,@('3+74=', '73-22=')
    ,@('23-0=', '87-60=')
    ,@('40-5=', '81-13=')
    ,@('33-25=', '14+56=')
    ,@('34-15=', '39+27=')
    ,@('6+45=', '57-46=')
    ,@('95-63=', '3+88=')
    ,@('41-39=', '95-25=')
    ,@('32-12=', '54-36=')
    ,@('21-5=', '83-38=')
    ,@('56-7=', '66-25=')
    ,@('99-50=', '6+15=')
    ,@('33-3=', '26+64=')
    ,@('79-16=', '14+65=')
    ,@('66-40=', '29+41=')
    ,@('76-6=', '11-1=')
    ,@('47-37=', '68-12=')
    ,@('37+7=', '6+1=')
    ,@('9+18=', '78-36=')
    ,@('46+51=', '76-2=')
    ,@('51-29=', '40+10=')
    ,@('3+35=', '49-39=')
    ,@('92-26=', '46+38=')
    ,@('81+8=', '49-45=')
    ,@('92+4=', '58+13=')
    ,@('39+1=', '54-22=')
    ,@('16+74=', '66-44=')
    ,@('25+70=', '4+23=')
    ,@('62+32=', '90-88=')
    ,@('44-17=', '96-26=')
    ,@('78-45=', '54+3=')
    ,@('23+8=', '41-8=')
    ,@('2+70=', '45+12=')
    ,@('5+51=', '25-18=')
    ,@('18-11=', '93-25=')
    ,@('86-79=', '32+12=')
    ,@('75+24=', '86-18=')
    ,@('31+10=', '30+12=')
    ,@('23+40=', '63-25=')
    ,@('51+27=', '10+45=')
    ,@('85-53=', '61-58=')
    ,@('3+19=', '67+10=')
    ,@('74+10=', '83-14=')
    ,@('25-6=', '32+57=')
    ,@('21+13=', '63-23=')
    ,@('71-43=', '84+3=')
    ,@('56+24=', '8-0=')
    ,@('19+62=', '15+66=')
    ,@('18+73=', '25-18=')
    ,@('88+9=', '27+59=')
    ,@('4+61=', '8+0=')
    ,@('29+37=', '47+50=')
    ,@('43-17=', '86+0=')
    ,@('22+31=', '66-6=')
    ,@('53-23=', '81-2=')
    ,@('3+54=', '42+48=')
    ,@('68-29=', '9+13=')
    ,@('34+35=', '27+47=')
    ,@('57+31=', '19+57=')
    ,@('3+14=', '0+26=')
    ,@('39-15=', '51+35=')
    ,@('23+48=', '97-82=')
    ,@('31+36=', '47-44=')
    ,@('69+18=', '96-43=')
    ,@('28+45=', '86-71=')
    ,@('78-3=', '17+48=')
    ,@('83-59=', '17+6=')
    ,@('29+50=', '51-51=')
    ,@('58-31=', '45-2=')
    ,@('76+10=', '88-35=')
    ,@('55-53=', '49+11=')
    ,@('16-0=', '85-49=')
    ,@('12+30=', '77+6=')
    ,@('10+60=', '5+32=')
)

$count = 0
foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $found = $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
    if ($found) { $count = $count + 1 } else { Write-Output "NOT FOUND: $old" }
}
Write-Output "Replaced $count of $($replacements.Count)"
